# Data cleaning and preparation fix
# - Rename the raw data sheet (drop the stray "_2" suffix picked up during
#   an earlier save).
# - Scroll the sheet view down a few rows (the view had been left scrolled
#   to row 4 when the workbook was last saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "raw_home_loan_approval"

$win = $excel.Windows.Item(1)
$win.ScrollRow = 4
$win.ScrollColumn = 1
